# "Invoice Entry with AutoNumber"
# Adds a new "Sale Order Entry" task row (row 13) to the Tasks sheet and
# updates the active selection / scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Add the new task row (row 13) -----------------------------------
# Copy formatting from the matching cells in row 12 (row 10 for the
# "Status" column, which is unstyled on most rows) so the new row picks
# up the same cell styles already used elsewhere in the sheet, then fill
# in the values for the new "Sale Order Entry" task.

$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)

$ws.Range("F10").Copy()
$ws.Range("F13").PasteSpecial(-4122)

$ws.Range("G12").Copy()
$ws.Range("G13").PasteSpecial(-4122)

$ws.Range("H12").Copy()
$ws.Range("H13").PasteSpecial(-4122)

$ws.Range("J12").Copy()
$ws.Range("J13").PasteSpecial(-4122)

$ws.Range("A13").Value = "Sale Order Entry"
$ws.Range("C13").Value = "10/31/2014"
$ws.Range("F13").Value = "Completed & Revised"
$ws.Range("G13").Value = "Done"
$ws.Range("H13").Value = "Done"
$ws.Range("J13").Value = "Done"

$excel.CutCopyMode = $false

# --- Update the view / selection --------------------------------------
$ws.Range("C10").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
